# Auto-generated Excel COM-interop edit script.
# Applies the numeric corrections described by the commit diff to the
# "Maduin_Profits" workbook, which is split across 8 per-job worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Each worksheet has identical
# columns: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# row 17
$ws.Range("H17").Value = 318
$ws.Range("J17").Value = 318
$ws.Range("L17").Value = 954
$ws.Range("N17").Value = -1290

# row 70
$ws.Range("H70").Value = 2999.5
$ws.Range("J70").Value = 2999
$ws.Range("L70").Value = 8997
$ws.Range("N70").Value = -9537

# row 73
$ws.Range("H73").Value = 2999.5
$ws.Range("J73").Value = 2999
$ws.Range("L73").Value = 8997
$ws.Range("N73").Value = -10869

# row 113
$ws.Range("H113").Value = 15199.4
$ws.Range("I113").Value = 20332.334
$ws.Range("K113").Value = 20332.334
$ws.Range("M113").Value = -17078.334

# row 132
$ws.Range("H132").Value = 1377.6364
$ws.Range("I132").Value = 1377.6364
$ws.Range("K132").Value = 4132.9092
$ws.Range("M132").Value = -1602.9092

# row 138
$ws.Range("H138").Value = 6411.9
$ws.Range("I138").Value = 3885.818
$ws.Range("J138").Value = 9499.333000000001
$ws.Range("K138").Value = 11657.454
$ws.Range("L138").Value = 28497.999
$ws.Range("M138").Value = -6517.454000000002
$ws.Range("N138").Value = -38777.999

$ws = $wb.Worksheets.Item("ARM")
# row 102
$ws.Range("H102").Value = 1320.091
$ws.Range("I102").Value = 1242.1
$ws.Range("K102").Value = 1242.1
$ws.Range("M102").Value = 379.9000000000001

# row 119
$ws.Range("H119").Value = 50348.5
$ws.Range("J119").Value = 50348.5
$ws.Range("L119").Value = 50348.5
$ws.Range("N119").Value = -60024.5

# row 132
$ws.Range("H132").Value = 2598
$ws.Range("I132").Value = 2598
$ws.Range("K132").Value = 7794
$ws.Range("M132").Value = -5264

$ws = $wb.Worksheets.Item("BSM")
# row 22
$ws.Range("H22").Value = 346.875
$ws.Range("I22").Value = 343.7143
$ws.Range("J22").Value = 369
$ws.Range("K22").Value = 343.7143
$ws.Range("L22").Value = 369
$ws.Range("M22").Value = -170.7143
$ws.Range("N22").Value = -715

# row 86
$ws.Range("H86").Value = 7081
$ws.Range("I86").Value = 7656
$ws.Range("J86").Value = 1906
$ws.Range("K86").Value = 7656
$ws.Range("L86").Value = 1906
$ws.Range("M86").Value = -6533
$ws.Range("N86").Value = -4152

# row 89
$ws.Range("H89").Value = 7081
$ws.Range("I89").Value = 7656
$ws.Range("J89").Value = 1906
$ws.Range("K89").Value = 38280
$ws.Range("L89").Value = 9530
$ws.Range("M89").Value = -32664
$ws.Range("N89").Value = -20762

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 949.75
$ws.Range("I16").Value = 933
$ws.Range("K16").Value = 933
$ws.Range("M16").Value = -646

# row 31
$ws.Range("H31").Value = 2991.7058
$ws.Range("I31").Value = 1925.8334
$ws.Range("K31").Value = 1925.8334
$ws.Range("M31").Value = -1630.8334

# row 34
$ws.Range("H34").Value = 2991.7058
$ws.Range("I34").Value = 1925.8334
$ws.Range("K34").Value = 1925.8334
$ws.Range("M34").Value = -1723.8334

# row 50
$ws.Range("H50").Value = 21581.572
$ws.Range("I50").Value = 15541.5
$ws.Range("J50").Value = 23997.6
$ws.Range("K50").Value = 15541.5
$ws.Range("L50").Value = 23997.6
$ws.Range("M50").Value = -14916.5
$ws.Range("N50").Value = -25247.6

# row 51
$ws.Range("H51").Value = 20033
$ws.Range("J51").Value = 20033
$ws.Range("L51").Value = 20033
$ws.Range("N51").Value = -21505

# row 60
$ws.Range("H60").Value = 20000
$ws.Range("J60").Value = 20000
$ws.Range("L60").Value = 20000
$ws.Range("N60").Value = -21022

# row 61
$ws.Range("H61").Value = 20033
$ws.Range("J61").Value = 20033
$ws.Range("L61").Value = 20033
$ws.Range("N61").Value = -20729

# row 69
$ws.Range("H69").Value = 3276
$ws.Range("I69").Value = 3276
$ws.Range("K69").Value = 3276
$ws.Range("M69").Value = -2527

# row 72
$ws.Range("H72").Value = 3276
$ws.Range("I72").Value = 3276
$ws.Range("K72").Value = 9828
$ws.Range("M72").Value = -6084

# row 99
$ws.Range("H99").Value = 4853.125
$ws.Range("I99").Value = 4689.2856
$ws.Range("K99").Value = 4689.2856
$ws.Range("M99").Value = -3191.2856

# row 107
$ws.Range("H107").Value = 500.63635
$ws.Range("I107").Value = 541
$ws.Range("K107").Value = 541
$ws.Range("M107").Value = 1379

# row 113
$ws.Range("H113").Value = 949.75
$ws.Range("I113").Value = 933
$ws.Range("K113").Value = 933
$ws.Range("M113").Value = 1237

# row 126
$ws.Range("H126").Value = 4853.125
$ws.Range("I126").Value = 4689.2856
$ws.Range("K126").Value = 14067.8568
$ws.Range("M126").Value = -11597.8568

# row 134
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 70.35294
$ws.Range("J2").Value = 125.28571
$ws.Range("L2").Value = 751.71426
$ws.Range("N2").Value = -977.71426

# row 22
$ws.Range("H22").Value = 129998.75
$ws.Range("J22").Value = 15000
$ws.Range("L22").Value = 45000
$ws.Range("N22").Value = -45338

# row 27
$ws.Range("H27").Value = 129998.75
$ws.Range("J27").Value = 15000
$ws.Range("L27").Value = 45000
$ws.Range("N27").Value = -45204

# row 38
$ws.Range("H38").Value = 744.44446
$ws.Range("I38").Value = 903.3333
$ws.Range("J38").Value = 426.66666
$ws.Range("K38").Value = 2709.9999
$ws.Range("L38").Value = 1279.99998
$ws.Range("M38").Value = -2362.9999
$ws.Range("N38").Value = -1973.99998

# row 68
$ws.Range("H68").Value = 3572.3845
$ws.Range("I68").Value = 2720.5715
$ws.Range("J68").Value = 4566.1665
$ws.Range("K68").Value = 8161.7145
$ws.Range("L68").Value = 13698.4995
$ws.Range("M68").Value = -7350.7145
$ws.Range("N68").Value = -15320.4995

# row 70
$ws.Range("H70").Value = 3450
$ws.Range("I70").Value = 3450
$ws.Range("K70").Value = 10350
$ws.Range("M70").Value = -10035

# row 71
$ws.Range("H71").Value = 3572.3845
$ws.Range("I71").Value = 2720.5715
$ws.Range("J71").Value = 4566.1665
$ws.Range("K71").Value = 24485.1435
$ws.Range("L71").Value = 41095.4985
$ws.Range("M71").Value = -20429.1435
$ws.Range("N71").Value = -49207.4985

# row 73
$ws.Range("H73").Value = 3450
$ws.Range("I73").Value = 3450
$ws.Range("K73").Value = 10350
$ws.Range("M73").Value = -9258

# row 108
$ws.Range("H108").Value = 320.2
$ws.Range("I108").Value = 320.2
$ws.Range("K108").Value = 960.5999999999999
$ws.Range("M108").Value = 1919.4

# row 109
$ws.Range("H109").Value = 3587.125
$ws.Range("I109").Value = 3049.25
$ws.Range("J109").Value = 4125
$ws.Range("K109").Value = 9147.75
$ws.Range("L109").Value = 12375
$ws.Range("M109").Value = -8107.75
$ws.Range("N109").Value = -14455

# row 112
$ws.Range("H112").Value = 44124.25
$ws.Range("I112").Value = 8000
$ws.Range("K112").Value = 24000
$ws.Range("M112").Value = -22892

# row 131
$ws.Range("H131").Value = 976
$ws.Range("J131").Value = 990
$ws.Range("L131").Value = 2970
$ws.Range("N131").Value = -13050

# row 140
$ws.Range("H140").Value = 799.2857
$ws.Range("I140").Value = 799.2857
$ws.Range("K140").Value = 2397.8571
$ws.Range("M140").Value = 2782.1429

$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 3237.3333
$ws.Range("I102").Value = 3512
$ws.Range("J102").Value = 3100
$ws.Range("K102").Value = 3512
$ws.Range("L102").Value = 3100
$ws.Range("M102").Value = -1890
$ws.Range("N102").Value = -6344

$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 6098.857
$ws.Range("I16").Value = 6331.3335
$ws.Range("K16").Value = 6331.3335
$ws.Range("M16").Value = -6161.3335

# row 20
$ws.Range("H20").Value = 206
$ws.Range("J20").Value = 206
$ws.Range("L20").Value = 206
$ws.Range("N20").Value = -658

# row 40
$ws.Range("H40").Value = 2272.818
$ws.Range("I40").Value = 2427.5715
$ws.Range("J40").Value = 2002
$ws.Range("K40").Value = 2427.5715
$ws.Range("L40").Value = 2002
$ws.Range("M40").Value = -2291.5715
$ws.Range("N40").Value = -2274

# row 46
$ws.Range("H46").Value = 4374.75
$ws.Range("I46").Value = 2999
$ws.Range("J46").Value = 4571.2856
$ws.Range("K46").Value = 2999
$ws.Range("L46").Value = 4571.2856
$ws.Range("M46").Value = -2811
$ws.Range("N46").Value = -4947.2856

$ws = $wb.Worksheets.Item("WVR")
# row 126
$ws.Range("H126").Value = 1949
$ws.Range("I126").Value = 1931.6666
$ws.Range("K126").Value = 5794.9998
$ws.Range("M126").Value = -3324.9998

# row 140
$ws.Range("H140").Value = 58749
$ws.Range("I140").Value = 79997
$ws.Range("K140").Value = 79997
$ws.Range("M140").Value = -74817
